$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clean up / correct the Genotype (column B) values ---
# These rows were mislabeled before Lizzie's cleanup pass; fix them up.
$ws.Range("B3").Value = "SAPAP3(WT)"
$ws.Range("B5").Value = "SAPAP3(WT)"
$ws.Range("B6").Value = "SAPAP3(WT)"
$ws.Range("B7").Value = "SAPAP3(WT)"
$ws.Range("B8").Value = "SAPAP3(WT)"
$ws.Range("B11").Value = "SAPAP3(WT)"
$ws.Range("B14").Value = "SAPAP3(WT)"
$ws.Range("B16").Value = "SAPAP3(KO)"
$ws.Range("B17").Value = "SAPAP3(KO)"
$ws.Range("B20").Value = "SAPAP3(KO)"
$ws.Range("B22").Value = "SAPAP3(KO)"
$ws.Range("B23").Value = "SAPAP3(KO)"
$ws.Range("B27").Value = "SAPAP3(KO)"
$ws.Range("B28").Value = "SAPAP3(KO)"

# --- Normalize the stray "applyFill" formatting left on some column B cells ---
# (visually identical - fill pattern is "none" either way - just tidies the style table)
$ws.Range("B9:B13").Interior.Pattern = -4142
$ws.Range("B23:B25").Interior.Pattern = -4142

# --- Add a 4th (currently blank/placeholder) column, matching the header & body formats ---
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("D2:D28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the saved selection / view ---
$ws.Range("G8").Select() | Out-Null
